$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $value) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $value
    $cellRange.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '61.354.61'
Set-TextValue $ws.Range("E2") '  +0.91%  '

# Row 3
Set-TextValue $ws.Range("D3") '2.377.60'
Set-TextValue $ws.Range("E3") '  +0.79%  '

# Row 4
Set-TextValue $ws.Range("E4") '  +0.02%  '

# Row 5
Set-TextValue $ws.Range("D5") '552.83'
Set-TextValue $ws.Range("E5") '  +2.59%  '

# Row 6
Set-TextValue $ws.Range("D6") '139.58'
Set-TextValue $ws.Range("E6") '  +1.42%  '

# Row 7
Set-TextValue $ws.Range("E7") '  -0.04%  '

# Row 8
Set-TextValue $ws.Range("E8") '  +0.11%  '

# Row 9
Set-TextValue $ws.Range("D9") '2.378.58'
Set-TextValue $ws.Range("E9") '  +0.83%  '

# Row 10
Set-TextValue $ws.Range("E10") '  +3.13%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.157'
Set-TextValue $ws.Range("E11") '  +2.20%  '

# Row 12
Set-TextValue $ws.Range("E12") '  +2.43%  '

# Row 13
Set-TextValue $ws.Range("E13") '  +3.25%  '

# Row 14
Set-TextValue $ws.Range("D14") '25.53'
Set-TextValue $ws.Range("E14") '  +3.17%  '

# Row 15
Set-TextValue $ws.Range("D15") '0.0000172'
Set-TextValue $ws.Range("E15") '  +6.72%  '

# Row 16
Set-TextValue $ws.Range("D16") '2.807.18'
Set-TextValue $ws.Range("E16") '  +0.79%  '

# Row 17
Set-TextValue $ws.Range("D17") '61.355.22'
Set-TextValue $ws.Range("E17") '  +1.48%  '

# Row 18
Set-TextValue $ws.Range("D18") '2.381.30'
Set-TextValue $ws.Range("E18") '  +0.82%  '

# Row 19
Set-TextValue $ws.Range("D19") '10.94'
Set-TextValue $ws.Range("E19") '  +3.39%  '

# Row 20
Set-TextValue $ws.Range("D20") '4.15'
Set-TextValue $ws.Range("E20") '  +2.47%  '

# Row 21
Set-TextValue $ws.Range("D21") '320.71'
Set-TextValue $ws.Range("E21") '  +1.83%  '

# Row 22
Set-TextValue $ws.Range("D22") '6.69'
Set-TextValue $ws.Range("E22") '  +1.78%  '

# Row 24
Set-TextValue $ws.Range("B24") 'Litecoin'
Set-TextValue $ws.Range("C24") 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range("D24") '64.28'
Set-TextValue $ws.Range("E24") '  +1.68%  '

# Row 25
Set-TextValue $ws.Range("B25") 'SuiNetwork'
Set-TextValue $ws.Range("C25") 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue $ws.Range("D25") '1.74'
Set-TextValue $ws.Range("E25") '  -8.10%  '

# Row 26
Set-TextValue $ws.Range("E26") '  +5.01%  '

# Row 27
Set-TextValue $ws.Range("D27") '0.999'
Set-TextValue $ws.Range("E27") '  -0.09%  '

# Row 28
Set-TextValue $ws.Range("D28") '2.494.60'
Set-TextValue $ws.Range("E28") '  +0.73%  '

# Row 29
Set-TextValue $ws.Range("D29") '8.16'
Set-TextValue $ws.Range("E29") '  +2.77%  '

# Row 30
Set-TextValue $ws.Range("D30") '519.49'
Set-TextValue $ws.Range("E30") '  +3.14%  '

# Row 31
Set-TextValue $ws.Range("D31") '0.0₃0902'
Set-TextValue $ws.Range("E31") '  +0.65%  '

# Row 32
Set-TextValue $ws.Range("E32") '  +0.52%  '

# Row 33
Set-TextValue $ws.Range("E33") '  +2.92%  '

# Row 34
Set-TextValue $ws.Range("D34") '1.83'
Set-TextValue $ws.Range("E34") '  +3.39%  '

# Row 35
Set-TextValue $ws.Range("E35") '  -0.70%  '

# Row 36
Set-TextValue $ws.Range("E36") '  +0.00%  '

# Row 37
Set-TextValue $ws.Range("D37") '5.52'

# Row 38
Set-TextValue $ws.Range("D38") '4.68'
Set-TextValue $ws.Range("E38") '  +2.85%  '

# Row 39
Set-TextValue $ws.Range("D39") '1.88'
Set-TextValue $ws.Range("E39") '  +6.10%  '

# Row 40
Set-TextValue $ws.Range("E40") '  +1.78%  '

# Row 41
Set-TextValue $ws.Range("D41") '18.49'

# Row 42
Set-TextValue $ws.Range("D42") '146.54'
Set-TextValue $ws.Range("E42") '  +6.22%  '

# Row 43
Set-TextValue $ws.Range("E43") '  +0.01%  '

# Row 44
Set-TextValue $ws.Range("E44") '  +2.93%  '

# Row 45
Set-TextValue $ws.Range("D45") '147.87'
Set-TextValue $ws.Range("E45") '  +6.81%  '

# Row 46
Set-TextValue $ws.Range("E46") '  +2.07%  '

# Row 47
Set-TextValue $ws.Range("D47") '3.60'
Set-TextValue $ws.Range("E47") '  +2.61%  '

# Row 48
Set-TextValue $ws.Range("D48") '0.0522'
Set-TextValue $ws.Range("E48") '  +2.41%  '

# Row 49
Set-TextValue $ws.Range("D49") '19.72'
Set-TextValue $ws.Range("E49") '  +1.50%  '

# Row 50
Set-TextValue $ws.Range("D50") '0.581'
Set-TextValue $ws.Range("E50") '  +2.07%  '

# Row 51
Set-TextValue $ws.Range("E51") '  +1.31%  '
